$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "purpose" column (E2:E25) to the new value "fullRNASEQ"
$ws.Range("E2:E25").Value = "fullRNASEQ"

# Update the sheet view: scroll so row 12 is the top-left visible row,
# and select E24:E25 with E24 as the active cell
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("E24:E25").Select()
$excel.ActiveWindow.RangeSelection.Item(1).Activate()
